# New feed-log runs to append under the existing data (run_id, rss_url_id,
# date, response, item_count) starting right after the current last row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3, 1, "2024-06-14 19:10:20", 200, 0),
    @(4, 2, "2024-06-14 19:10:20", 200, 0)
)

$lastRow = 3
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}
